$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto price/volume values (GitHub Actions data refresh)
$ws.Range("D2").Value = "'28.269.50"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -0.78%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.808.40"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -0.92%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'0.9979"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  -0.64%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'312.67"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -1.21%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'0.9973"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -0.61%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.5167"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -0.20%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.3978"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +2.94%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.07880"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -6.24%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D11").Value = "'41.10"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -2.06%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'6.346"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -1.07%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.9974"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -0.63%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'20.45"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -3.55%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'7.328"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -2.34%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'1.793.34"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -1.45%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'92.41"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -1.92%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'0.00001082"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -4.21%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'0.06564"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -1.11%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'0.9966"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -0.65%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'17.32"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -2.45%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'5.997"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -1.25%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'28.330.91"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -0.73%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'11.13"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -2.41%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'2.223"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -2.98%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'160.41"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +0.49%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'20.59"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -2.89%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'2.009.97"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -1.07%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Value = "'  -0.16%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'127.76"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +1.44%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'0.1086"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -0.82%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'1.049"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -4.68%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'5.582"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Value = "'3.654"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -0.59%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'0.07156"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -6.88%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'9.069"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +3.56%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'0.02330"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -2.16%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'0.2152"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -3.58%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = "'  +1.04%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'5.059"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -4.36%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.6203"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -3.32%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.9964"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -0.54%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'1.154"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -3.58%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'13.22"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -3.01%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'1.326"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -5.29%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.5975"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -2.86%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'3.747"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -1.27%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'125.55"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -1.82%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'1.213"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +0.58%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'1.938"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -3.05%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.06869"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -1.81%  "
$ws.Range("E51").Style = "Normal"
